# "Generate Report for Handoff"
# Update the localization-status report for file "b.md": it has now been
# handed off for localization (Ready for handoff), and a new handoff XLIFF
# was produced for it in both zh-cn and de-de, which is newer than the
# already-handed-back translation, hence the mismatch/error note.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/360cfff5ea70736d7978372665bcb73c21c7fab0/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7c313e230dfe7e6b954a39aac0514ca10d016eeb/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 is b.md
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-01 18:44:23"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is b.md
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces text (not boolean) storage for "False"; reset
# the style afterwards so no stray quote-prefix formatting is left behind.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-01 18:44:18"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: row 3 is b.md
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-01 18:44:23"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
